# Family tree data cleanup:
#  - One row had a bogus/placeholder parent id (99999) in "Parent 1" -
#    clear it out (change it to null/blank).
#  - Re-sort the data table by "Parse ID" (column D) ascending instead of
#    by "Last Name" (column B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the bogus parent id (row 19, column E = "Parent 1") before sorting
# so it travels with its row regardless of row order.
$ws.Range("E19").Value = ""

# Re-sort Table2 by column D ("Parse ID") ascending, replacing the
# previous sort-by-Last-Name.
$lo = $ws.ListObjects.Item("Table2")
$sort = $lo.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("D1:D151"))
$sort.Header = 1
$sort.Apply()

# Leave the selection where the user ended up after the re-sort/review.
$ws.Range("E3").Select()
